$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in cell B2 (this introduces the new shared string "${movDescID}")
$ws.Range("B2").Value = "`${movDescID}"

# Move the active selection from G2 to B3
$ws.Range("B3").Select()
